# "output folder feature - extra context can be saved - repeated analyses
# are saved seperately"
#
# Formats the data rows of the Ratios sheet (everything below the header
# row) with an Arial font, centered horizontally, while leaving the
# existing alternating-row green banding fill untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows run from row 2 (first repeated-analysis row) through the last
# used row; columns A:Y (Lab # plus the 24 ratio/error columns).
$lastRow = $ws.UsedRange.Rows.Count
$dataRange = $ws.Range("A2:Y" + $lastRow)

$dataRange.Font.Name = "Arial"
$dataRange.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
